$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 151.75
$ws.Range("I2").Value = 30.571428
$ws.Range("K2").Value = 30.571428
$ws.Range("M2").Value = 82.428572
$ws.Range("H11").Value = 3243.5833
$ws.Range("I11").Value = 3243.5833
$ws.Range("K11").Value = 3243.5833
$ws.Range("M11").Value = -3103.5833
$ws.Range("H53").Value = 5463.3335
$ws.Range("I53").Value = 6242
$ws.Range("J53").Value = 4907.143
$ws.Range("K53").Value = 6242
$ws.Range("L53").Value = 4907.143
$ws.Range("M53").Value = -5605
$ws.Range("N53").Value = -6181.143
$ws.Range("H62").Value = 3732.5557
$ws.Range("I62").Value = 1847.8334
$ws.Range("K62").Value = 1847.8334
$ws.Range("M62").Value = -1223.8334
$ws.Range("H65").Value = 3732.5557
$ws.Range("I65").Value = 1847.8334
$ws.Range("K65").Value = 9239.166999999999
$ws.Range("M65").Value = -6119.166999999999
$ws.Range("H86").Value = 154324180
$ws.Range("I86").Value = 444445250
$ws.Range("J86").Value = 9263641
$ws.Range("K86").Value = 444445250
$ws.Range("L86").Value = 9263641
$ws.Range("M86").Value = -444444127
$ws.Range("N86").Value = -9265887
$ws.Range("H89").Value = 154324180
$ws.Range("I89").Value = 444445250
$ws.Range("J89").Value = 9263641
$ws.Range("K89").Value = 2222226250
$ws.Range("L89").Value = 46318205
$ws.Range("M89").Value = -2222220634
$ws.Range("N89").Value = -46329437
$ws.Range("H95").Value = 53947
$ws.Range("J95").Value = 53947
$ws.Range("L95").Value = 53947
$ws.Range("N95").Value = -59439
$ws.Range("H98").Value = 1976.0588
$ws.Range("I98").Value = 2019.8372
$ws.Range("K98").Value = 2019.8372
$ws.Range("M98").Value = -521.8371999999999
$ws.Range("H112").Value = 8971.92
$ws.Range("I112").Value = 2327.5
$ws.Range("J112").Value = 9549.695
$ws.Range("K112").Value = 6982.5
$ws.Range("L112").Value = 28649.085
$ws.Range("M112").Value = -5874.5
$ws.Range("N112").Value = -30865.085
$ws.Range("H122").Value = 1976.0588
$ws.Range("I122").Value = 2019.8372
$ws.Range("K122").Value = 6059.5116
$ws.Range("M122").Value = -3609.5116
$ws.Range("H135").Value = 182526.69
$ws.Range("I135").Value = 227971.36
$ws.Range("J135").Value = 748
$ws.Range("K135").Value = 2051742.24
$ws.Range("L135").Value = 6732
$ws.Range("M135").Value = -2049207.24
$ws.Range("N135").Value = -11802
$ws.Range("H137").Value = 1125
$ws.Range("I137").Value = 1125
$ws.Range("K137").Value = 3375
$ws.Range("M137").Value = -825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2454875.8
$ws.Range("I32").Value = 2607180.2
$ws.Range("K32").Value = 2607180.2
$ws.Range("M32").Value = -2606893.2
$ws.Range("H61").Value = 9752.5
$ws.Range("I61").Value = 3632.7273
$ws.Range("J61").Value = 15872.272
$ws.Range("K61").Value = 3632.7273
$ws.Range("L61").Value = 15872.272
$ws.Range("M61").Value = -3420.7273
$ws.Range("N61").Value = -16296.272
$ws.Range("H74").Value = 24654.553
$ws.Range("I74").Value = 36472.414
$ws.Range("J74").Value = 5614.6665
$ws.Range("K74").Value = 36472.414
$ws.Range("L74").Value = 5614.6665
$ws.Range("M74").Value = -35598.414
$ws.Range("N74").Value = -7362.6665
$ws.Range("H77").Value = 24654.553
$ws.Range("I77").Value = 36472.414
$ws.Range("J77").Value = 5614.6665
$ws.Range("K77").Value = 182362.07
$ws.Range("L77").Value = 28073.3325
$ws.Range("M77").Value = -177994.07
$ws.Range("N77").Value = -36809.3325
$ws.Range("H97").Value = 7590035.5
$ws.Range("I97").Value = 408.8
$ws.Range("K97").Value = 408.8
$ws.Range("M97").Value = 87.19999999999999
$ws.Range("H103").Value = 54802
$ws.Range("J103").Value = 54802
$ws.Range("L103").Value = 54802
$ws.Range("N103").Value = -57146
$ws.Range("H122").Value = 4398.96
$ws.Range("I122").Value = 3817.7273
$ws.Range("J122").Value = 8661.333000000001
$ws.Range("K122").Value = 11453.1819
$ws.Range("L122").Value = 25983.999
$ws.Range("M122").Value = -9003.1819
$ws.Range("N122").Value = -30883.999
$ws.Range("H132").Value = 930755.5600000001
$ws.Range("I132").Value = 1354737
$ws.Range("J132").Value = 7972.4707
$ws.Range("K132").Value = 4064211
$ws.Range("L132").Value = 23917.4121
$ws.Range("M132").Value = -4061681
$ws.Range("N132").Value = -28977.4121
$ws.Range("H136").Value = 9752.5
$ws.Range("I136").Value = 3632.7273
$ws.Range("J136").Value = 15872.272
$ws.Range("K136").Value = 10898.1819
$ws.Range("L136").Value = 47616.81600000001
$ws.Range("M136").Value = -8348.1819
$ws.Range("N136").Value = -52716.81600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4908.0527
$ws.Range("I134").Value = 1317.7084
$ws.Range("J134").Value = 11062.929
$ws.Range("K134").Value = 3953.1252
$ws.Range("L134").Value = 33188.787
$ws.Range("M134").Value = -1418.1252
$ws.Range("N134").Value = -38258.787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8550.714
$ws.Range("I31").Value = 2129
$ws.Range("J31").Value = 11119.4
$ws.Range("K31").Value = 2129
$ws.Range("L31").Value = 11119.4
$ws.Range("M31").Value = -1834
$ws.Range("N31").Value = -11709.4
$ws.Range("H34").Value = 8550.714
$ws.Range("I34").Value = 2129
$ws.Range("J34").Value = 11119.4
$ws.Range("K34").Value = 2129
$ws.Range("L34").Value = 11119.4
$ws.Range("M34").Value = -1927
$ws.Range("N34").Value = -11523.4
$ws.Range("H134").Value = 12220.091
$ws.Range("I134").Value = 3012
$ws.Range("J134").Value = 13140.9
$ws.Range("K134").Value = 9036
$ws.Range("L134").Value = 39422.7
$ws.Range("M134").Value = -6501
$ws.Range("N134").Value = -44492.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 11906601
$ws.Range("I14").Value = 11906601
$ws.Range("K14").Value = 35719803
$ws.Range("M14").Value = -35719630
$ws.Range("H39").Value = 9266.5
$ws.Range("J39").Value = 9652
$ws.Range("L39").Value = 28956
$ws.Range("N39").Value = -29544
$ws.Range("H68").Value = 2623.4546
$ws.Range("J68").Value = 2926.3125
$ws.Range("L68").Value = 8778.9375
$ws.Range("N68").Value = -10400.9375
$ws.Range("H71").Value = 2623.4546
$ws.Range("J71").Value = 2926.3125
$ws.Range("L71").Value = 26336.8125
$ws.Range("N71").Value = -34448.8125
$ws.Range("H107").Value = 1414.027
$ws.Range("J107").Value = 1780
$ws.Range("L107").Value = 5340
$ws.Range("N107").Value = -9180
$ws.Range("I129").Value = 491.9091
$ws.Range("K129").Value = 1475.7273
$ws.Range("M129").Value = 3524.2727
$ws.Range("H134").Value = 171465.67
$ws.Range("I134").Value = 182998.92
$ws.Range("K134").Value = 548996.76
$ws.Range("M134").Value = -543926.76
$ws.Range("H136").Value = 2030
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 203582.2
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 203582.2
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H122").Value = 25668966
$ws.Range("I122").Value = 34519216
$ws.Range("K122").Value = 103557648
$ws.Range("M122").Value = -103555198
$ws.Range("H126").Value = 6500
$ws.Range("I126").Value = 2950
$ws.Range("K126").Value = 8850
$ws.Range("M126").Value = -6380
$ws.Range("H132").Value = 3123.6099
$ws.Range("I132").Value = 1767.4814
$ws.Range("J132").Value = 5739
$ws.Range("K132").Value = 5302.4442
$ws.Range("L132").Value = 17217
$ws.Range("M132").Value = -2772.4442
$ws.Range("N132").Value = -22277

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2415
$ws.Range("I68").Value = 2496.6667
$ws.Range("J68").Value = 2333.3333
$ws.Range("K68").Value = 2496.6667
$ws.Range("L68").Value = 2333.3333
$ws.Range("M68").Value = -1747.6667
$ws.Range("N68").Value = -3831.3333
$ws.Range("H69").Value = 42999.5
$ws.Range("J69").Value = 42999.5
$ws.Range("L69").Value = 42999.5
$ws.Range("N69").Value = -44621.5
$ws.Range("H71").Value = 2415
$ws.Range("I71").Value = 2496.6667
$ws.Range("J71").Value = 2333.3333
$ws.Range("K71").Value = 12483.3335
$ws.Range("L71").Value = 11666.6665
$ws.Range("M71").Value = -8739.333500000001
$ws.Range("N71").Value = -19154.6665
$ws.Range("H72").Value = 42999.5
$ws.Range("J72").Value = 42999.5
$ws.Range("L72").Value = 128998.5
$ws.Range("N72").Value = -137110.5
$ws.Range("H100").Value = 3961.4
$ws.Range("I100").Value = 3334.3333
$ws.Range("K100").Value = 3334.3333
$ws.Range("M100").Value = -2793.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2234.8076
$ws.Range("I126").Value = 1705.7646
$ws.Range("K126").Value = 5117.293799999999
$ws.Range("M126").Value = -2647.293799999999
$ws.Range("H135").Value = 174444
$ws.Range("J135").Value = 174444
$ws.Range("L135").Value = 174444
$ws.Range("N135").Value = -184584
$ws.Range("H136").Value = 26407.117
$ws.Range("I136").Value = 1145.3334
$ws.Range("K136").Value = 3436.0002
$ws.Range("M136").Value = -886.0001999999999
$ws.Range("H140").Value = 74375.25
$ws.Range("J140").Value = 74375.25
$ws.Range("L140").Value = 74375.25
$ws.Range("N140").Value = -84735.25
$ws.Range("H141").Value = 71211.75
$ws.Range("J141").Value = 84949
$ws.Range("L141").Value = 84949
$ws.Range("N141").Value = -95309
